$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("D2")
$cell.NumberFormat = "@"
$cell.Value = "41.185.13"
$cell.Style = "Normal"
$ws.Range("E2").Value = "  -1.58%  "

$cell = $ws.Range("D3")
$cell.NumberFormat = "@"
$cell.Value = "2.167.23"
$cell.Style = "Normal"
$ws.Range("E3").Value = "  -1.84%  "

$ws.Range("E4").Value = "  -0.06%  "

$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = "236.54"
$cell.Style = "Normal"
$ws.Range("E5").Value = "  -1.57%  "

$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = "0.605"
$cell.Style = "Normal"
$ws.Range("E6").Value = "  -2.71%  "

$cell = $ws.Range("D7")
$cell.NumberFormat = "@"
$cell.Value = "69.02"
$cell.Style = "Normal"
$ws.Range("E7").Value = "  -4.92%  "

$ws.Range("E8").Value = "  +0.00%  "

$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = "0.569"
$cell.Style = "Normal"
$ws.Range("E9").Value = "  -5.00%  "

$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = "39.01"
$cell.Style = "Normal"
$ws.Range("E10").Value = "  -7.50%  "

$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = "0.0911"
$cell.Style = "Normal"
$ws.Range("E11").Value = "  -3.66%  "

$cell = $ws.Range("D12")
$cell.NumberFormat = "@"
$cell.Value = "54.26"
$cell.Style = "Normal"
$ws.Range("E12").Value = "  -5.60%  "

$ws.Range("E13").Value = "  -2.20%  "

$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = "6.66"
$cell.Style = "Normal"
$ws.Range("E14").Value = "  -5.05%  "

$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = "2.497.98"
$cell.Style = "Normal"
$ws.Range("E15").Value = "  -1.75%  "

$cell = $ws.Range("D16")
$cell.NumberFormat = "@"
$cell.Value = "13.90"
$cell.Style = "Normal"
$ws.Range("E16").Value = "  -1.61%  "

$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = "2.163.01"
$cell.Style = "Normal"
$ws.Range("E17").Value = "  -2.77%  "

$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value = "0.788"
$cell.Style = "Normal"
$ws.Range("E18").Value = "  -4.93%  "

$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = "41.017.07"
$cell.Style = "Normal"
$ws.Range("E19").Value = "  -1.83%  "

$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = "0.0₃0993"
$cell.Style = "Normal"
$ws.Range("E20").Value = "  -7.42%  "

$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = "70.26"
$cell.Style = "Normal"
$ws.Range("E21").Value = "  -3.52%  "

$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = "5.83"
$cell.Style = "Normal"
$ws.Range("E22").Value = "  -4.23%  "

$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = "223.96"
$cell.Style = "Normal"
$ws.Range("E23").Value = "  -1.82%  "

$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = "9.30"
$cell.Style = "Normal"
$ws.Range("E24").Value = "  -9.58%  "

$ws.Range("B25").Value = "Dai"
$ws.Range("C25").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = "1.00"
$cell.Style = "Normal"
$ws.Range("E25").Value = "  -0.06%  "

$ws.Range("B26").Value = "ImmutableX"
$ws.Range("C26").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$cell = $ws.Range("D26")
$cell.NumberFormat = "@"
$cell.Value = "1.89"
$cell.Style = "Normal"
$ws.Range("E26").Value = "  -8.48%  "

$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = "10.64"
$cell.Style = "Normal"
$ws.Range("E27").Value = "  -8.02%  "

$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = "3.45"
$cell.Style = "Normal"
$ws.Range("E28").Value = "  -4.12%  "

$ws.Range("E29").Value = "  -2.38%  "

$cell = $ws.Range("D30")
$cell.NumberFormat = "@"
$cell.Value = "2.17"
$cell.Style = "Normal"
$ws.Range("E30").Value = "  +0.03%  "

$cell = $ws.Range("D31")
$cell.NumberFormat = "@"
$cell.Value = "167.55"
$cell.Style = "Normal"
$ws.Range("E31").Value = "  +0.59%  "

$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value = "19.75"
$cell.Style = "Normal"
$ws.Range("E32").Value = "  -3.52%  "

$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = "29.99"
$cell.Style = "Normal"
$ws.Range("E33").Value = "  +4.17%  "

$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = "0.0756"
$cell.Style = "Normal"
$ws.Range("E34").Value = "  -3.89%  "

$cell = $ws.Range("D35")
$cell.NumberFormat = "@"
$cell.Value = "5.05"
$cell.Style = "Normal"
$ws.Range("E35").Value = "  -10.96%  "

$ws.Range("E36").Value = "  -3.23%  "

$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value = "0.100"
$cell.Style = "Normal"
$ws.Range("E37").Value = "  -8.03%  "

$cell = $ws.Range("D38")
$cell.NumberFormat = "@"
$cell.Value = "4.03"
$cell.Style = "Normal"
$ws.Range("E38").Value = "  -4.56%  "

$cell = $ws.Range("D39")
$cell.NumberFormat = "@"
$cell.Value = "0.0279"
$cell.Style = "Normal"
$ws.Range("E39").Value = "  -5.97%  "

$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = "2.05"
$cell.Style = "Normal"
$ws.Range("E40").Value = "  -2.56%  "

$ws.Range("B41").Value = "Celestia"
$ws.Range("C41").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = "11.56"
$cell.Style = "Normal"
$ws.Range("E41").Value = "  -13.40%  "

$ws.Range("B42").Value = "THORChain"
$ws.Range("C42").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = "5.37"
$cell.Style = "Normal"
$ws.Range("E42").Value = "  -3.91%  "

$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = "58.46"
$cell.Style = "Normal"
$ws.Range("E43").Value = "  -10.76%  "

$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = "0.187"
$cell.Style = "Normal"
$ws.Range("E44").Value = "  -4.74%  "

$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = "8.21"
$cell.Style = "Normal"
$ws.Range("E45").Value = "  -4.70%  "

$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = "0.0966"
$cell.Style = "Normal"
$ws.Range("E46").Value = "  -3.52%  "

$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = "96.62"
$cell.Style = "Normal"
$ws.Range("E47").Value = "  -6.51%  "

$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = "1.07"
$cell.Style = "Normal"
$ws.Range("E48").Value = "  -3.32%  "

$ws.Range("E49").Value = "  -3.04%  "

$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = "2.18"
$cell.Style = "Normal"
$ws.Range("E50").Value = "  -7.64%  "

$ws.Range("E51").Value = "  -2.63%  "
